$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 66.85111099999999
$ws.Range("H2").Value = 200.553333
$ws.Range("I2").Value = 0.1215550702639512
$ws.Range("J2").Value = 0.1215550702639512
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 62.49926866666667
$ws.Range("N2").Value = 187.497806
$ws.Range("O2").Value = 0.364999568828264
$ws.Range("P2").Value = 0.364999568828264
$ws.Range("Q2").Value = 4178.145547054155
$ws.Range("R2").Value = 37603.3099234874
$ws.Range("S2").Value = 0.04436754823523154
$ws.Range("T2").Value = 0.04436754823523154

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 66.85111099999999
$ws.Range("H3").Value = 200.553333
$ws.Range("I3").Value = 0.1215550702639512
$ws.Range("J3").Value = 0.1215550702639512
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 43.717953
$ws.Range("N3").Value = 131.153859
$ws.Range("O3").Value = 0.2553155314530077
$ws.Range("P3").Value = 0.2553155314530077
$ws.Range("Q3").Value = 2922.593728695783
$ws.Range("R3").Value = 26303.34355826205
$ws.Range("S3").Value = 0.0310348973652484
$ws.Range("T3").Value = 0.0310348973652484

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 66.85111099999999
$ws.Range("H4").Value = 200.553333
$ws.Range("I4").Value = 0.1215550702639512
$ws.Range("J4").Value = 0.1215550702639512
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 45.87732766666667
$ws.Range("N4").Value = 137.631983
$ws.Range("O4").Value = 0.2679264121734788
$ws.Range("P4").Value = 0.2679264121734788
$ws.Range("Q4").Value = 3066.950324227704
$ws.Range("R4").Value = 27602.55291804934
$ws.Range("S4").Value = 0.03256781385731557
$ws.Range("T4").Value = 0.03256781385731557

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 66.85111099999999
$ws.Range("H5").Value = 200.553333
$ws.Range("I5").Value = 0.1215550702639512
$ws.Range("J5").Value = 0.1215550702639512
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.136526
$ws.Range("N5").Value = 57.409578
$ws.Range("O5").Value = 0.1117584875452494
$ws.Range("P5").Value = 0.1117584875452494
$ws.Range("Q5").Value = 1279.298023780386
$ws.Range("R5").Value = 11513.68221402347
$ws.Range("S5").Value = 0.01358481080615571
$ws.Range("T5").Value = 0.01358481080615571

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 99.86393233333332
$ws.Range("H6").Value = 299.591797
$ws.Range("I6").Value = 0.1815821327429069
$ws.Range("J6").Value = 0.1815821327429069
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 62.49926866666667
$ws.Range("N6").Value = 187.497806
$ws.Range("O6").Value = 0.364999568828264
$ws.Range("P6").Value = 0.364999568828264
$ws.Range("Q6").Value = 6241.42273701082
$ws.Range("R6").Value = 56172.80463309739
$ws.Range("S6").Value = 0.06627740015807762
$ws.Range("T6").Value = 0.06627740015807762

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 99.86393233333332
$ws.Range("H7").Value = 299.591797
$ws.Range("I7").Value = 0.1815821327429069
$ws.Range("J7").Value = 0.1815821327429069
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 43.717953
$ws.Range("N7").Value = 131.153859
$ws.Range("O7").Value = 0.2553155314530077
$ws.Range("P7").Value = 0.2553155314530077
$ws.Range("Q7").Value = 4365.846700143847
$ws.Range("R7").Value = 39292.62030129463
$ws.Range("S7").Value = 0.04636073872362587
$ws.Range("T7").Value = 0.04636073872362587

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 99.86393233333332
$ws.Range("H8").Value = 299.591797
$ws.Range("I8").Value = 0.1815821327429069
$ws.Range("J8").Value = 0.1815821327429069
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 45.87732766666667
$ws.Range("N8").Value = 137.631983
$ws.Range("O8").Value = 0.2679264121734788
$ws.Range("P8").Value = 0.2679264121734788
$ws.Range("Q8").Value = 4581.490345738162
$ws.Range("R8").Value = 41233.41311164346
$ws.Range("S8").Value = 0.04865064934061542
$ws.Range("T8").Value = 0.04865064934061542

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 99.86393233333332
$ws.Range("H9").Value = 299.591797
$ws.Range("I9").Value = 0.1815821327429069
$ws.Range("J9").Value = 0.1815821327429069
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.136526
$ws.Range("N9").Value = 57.409578
$ws.Range("O9").Value = 0.1117584875452494
$ws.Range("P9").Value = 0.1117584875452494
$ws.Range("Q9").Value = 1911.048737559074
$ws.Range("R9").Value = 17199.43863803166
$ws.Range("S9").Value = 0.02029334452058799
$ws.Range("T9").Value = 0.02029334452058799

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 368.4456226666667
$ws.Range("H10").Value = 1105.336868
$ws.Range("I10").Value = 0.6699429954379058
$ws.Range("J10").Value = 0.6699429954379058
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 62.49926866666667
$ws.Range("N10").Value = 187.497806
$ws.Range("O10").Value = 0.364999568828264
$ws.Range("P10").Value = 0.364999568828264
$ws.Range("Q10").Value = 23027.58196010129
$ws.Range("R10").Value = 207248.2376409117
$ws.Range("S10").Value = 0.2445289044743513
$ws.Range("T10").Value = 0.2445289044743513

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 368.4456226666667
$ws.Range("H11").Value = 1105.336868
$ws.Range("I11").Value = 0.6699429954379058
$ws.Range("J11").Value = 0.6699429954379058
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 43.717953
$ws.Range("N11").Value = 131.153859
$ws.Range("O11").Value = 0.2553155314530077
$ws.Range("P11").Value = 0.2553155314530077
$ws.Range("Q11").Value = 16107.68841479707
$ws.Range("R11").Value = 144969.1957331736
$ws.Range("S11").Value = 0.1710468519234488
$ws.Range("T11").Value = 0.1710468519234488

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 368.4456226666667
$ws.Range("H12").Value = 1105.336868
$ws.Range("I12").Value = 0.6699429954379058
$ws.Range("J12").Value = 0.6699429954379058
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 45.87732766666667
$ws.Range("N12").Value = 137.631983
$ws.Range("O12").Value = 0.2679264121734788
$ws.Range("P12").Value = 0.2679264121734788
$ws.Range("Q12").Value = 16903.3005584277
$ws.Range("R12").Value = 152129.7050258493
$ws.Range("S12").Value = 0.1794954231284314
$ws.Range("T12").Value = 0.1794954231284314

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 368.4456226666667
$ws.Range("H13").Value = 1105.336868
$ws.Range("I13").Value = 0.6699429954379058
$ws.Range("J13").Value = 0.6699429954379058
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.136526
$ws.Range("N13").Value = 57.409578
$ws.Range("O13").Value = 0.1117584875452494
$ws.Range("P13").Value = 0.1117584875452494
$ws.Range("Q13").Value = 7050.769237746857
$ws.Range("R13").Value = 63456.92313972171
$ws.Range("S13").Value = 0.07487181591167427
$ws.Range("T13").Value = 0.07487181591167427

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 14.80496566666667
$ws.Range("H14").Value = 44.414897
$ws.Range("I14").Value = 0.02691980155523597
$ws.Range("J14").Value = 0.02691980155523597
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 62.49926866666667
$ws.Range("N14").Value = 187.497806
$ws.Range("O14").Value = 0.364999568828264
$ws.Range("P14").Value = 0.364999568828264
$ws.Range("Q14").Value = 925.299526801776
$ws.Range("R14").Value = 8327.695741215985
$ws.Range("S14").Value = 0.009825715960603561
$ws.Range("T14").Value = 0.009825715960603559

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 14.80496566666667
$ws.Range("H15").Value = 44.414897
$ws.Range("I15").Value = 0.02691980155523597
$ws.Range("J15").Value = 0.02691980155523597
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 43.717953
$ws.Range("N15").Value = 131.153859
$ws.Range("O15").Value = 0.2553155314530077
$ws.Range("P15").Value = 0.2553155314530077
$ws.Range("Q15").Value = 647.242793181947
$ws.Range("R15").Value = 5825.185138637524
$ws.Range("S15").Value = 0.006873043440684574
$ws.Range("T15").Value = 0.006873043440684574

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 14.80496566666667
$ws.Range("H16").Value = 44.414897
$ws.Range("I16").Value = 0.02691980155523597
$ws.Range("J16").Value = 0.02691980155523597
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 45.87732766666667
$ws.Range("N16").Value = 137.631983
$ws.Range("O16").Value = 0.2679264121734788
$ws.Range("P16").Value = 0.2679264121734788
$ws.Range("Q16").Value = 679.212260983417
$ws.Range("R16").Value = 6112.910348850753
$ws.Range("S16").Value = 0.007212525847116408
$ws.Range("T16").Value = 0.007212525847116407

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 14.80496566666667
$ws.Range("H17").Value = 44.414897
$ws.Range("I17").Value = 0.02691980155523597
$ws.Range("J17").Value = 0.02691980155523597
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.136526
$ws.Range("N17").Value = 57.409578
$ws.Range("O17").Value = 0.1117584875452494
$ws.Range("P17").Value = 0.1117584875452494
$ws.Range("Q17").Value = 283.315610409274
$ws.Range("R17").Value = 2549.840493683466
$ws.Range("S17").Value = 0.003008516306831425
$ws.Range("T17").Value = 0.003008516306831425
